# repex for mutate function
# Activate the "Interview-CV95" sheet (was previously inactive; becomes
# the active / tabSelected sheet, matching the workbook's activeTab change),
# then fill the new column F ("repex") values for rows 2-31 and select
# the range F2:F31, leaving the newly written cells with the default
# (unstyled) format.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Interview-CV95")
$ws.Activate()

$values = @(14,14,14,14,14,14,17,17,17,17,17,17,14,14,14,11,11,11,11,11,11,10,10,10,9,13,11,13,12,12)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i]
}

# Newly entered cells keep the workbook's default (unstyled) format,
# unlike the pre-existing A:E columns which carry style index 1.
$ws.Range("F2:F31").Style = "Normal"

# Match the author's final selection in the worksheet.
[void]$ws.Range("F2:F31").Select()
